# Update countries & provincias Spain
# Applies the periodic COVID-19 data refresh to the "Pais" sheet:
#  - Updates the "Datos actualizados..." timestamp
#  - Updates case figures for a number of countries
#  - Two pairs of countries swap ranking position (and therefore rows),
#    since the sheet is kept sorted by total cases (column B) descending

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

function Set-Row {
    param($row, $b, $c, $d, $e, $f, $g, $h)
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# --- Timestamp (A1) ---------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 20 de Septiembre de 2020 a las 13:33"

# --- Plain data refreshes (country stays on the same row) -------------
# Row 4: Estados Unidos
Set-Row 4 6968918 1515 4223918 2541156 0 20 203844

# Row 16: Iran
Set-Row 16 422140 3097 359570 38269 0 183 24301

# Row 33: Rumania
Set-Row 33 112781 1231 89771 18575 0 33 4435

# Row 65: Ghana
Set-Row 65 46004 127 45153 554 0 0 297

# Row 71: Estado de Palestina
Set-Row 71 35686 683 23700 11724 0 9 262

# Row 79: Bosnia y Herzegovina
Set-Row 79 25428 211 17878 6787 0 5 763

# Row 87: Madagascar
Set-Row 87 16053 33 14646 1184 0 4 223

# Row 98: Malasia
Set-Row 98 10219 52 9355 734 0 0 130

# Row 145: Malta
Set-Row 145 2731 32 2047 664 0 1 20

# Row 182: Gibraltar
Set-Row 182 350 0 323 27 0 0 0

# --- Ranking swaps (two countries exchange rows + get refreshed data) -

# Rows 41/42: Paises Bajos <-> Oman
$ws.Range("A41").Value = "Oman"
Set-Row 41 93475 1722 85418 7211 0 28 846
$ws.Range("A42").Value = "Paises Bajos"
Set-Row 42 91934 0 0 0 0 0 6275

# Rows 204/205: Santa Lucia <-> Timor Oriental (data identical, names swap)
$ws.Range("A204").Value = "Timor Oriental"
Set-Row 204 27 0 26 1 0 0 0
$ws.Range("A205").Value = "Santa Lucia"
Set-Row 205 27 0 26 1 0 0 0

# Rows 214/215: Montserrat <-> Islas Malvinas
$ws.Range("A214").Value = "Islas Malvinas"
Set-Row 214 13 0 13 0 0 0 0
$ws.Range("A215").Value = "Montserrat"
Set-Row 215 13 0 12 0 0 0 1
